# Updated symbol list on Mon Dec 26 20:39:43 UTC 2022 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Ensure the cell keeps a text representation (matches source inlineStr cells)
    # even when the string looks numeric, so Excel doesn't silently convert it
    # to a numeric cell type.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Column D price updates (rows unaffected by the row-shift below) ---
Set-TextValue $ws.Range("D2")  "242.59"
Set-TextValue $ws.Range("D4")  "5.413"
Set-TextValue $ws.Range("D6")  "3.434"
Set-TextValue $ws.Range("D7")  "6.525"
Set-TextValue $ws.Range("D8")  "0.8085"
Set-TextValue $ws.Range("D9")  "0.9273"
Set-TextValue $ws.Range("D10") "0.1418"
Set-TextValue $ws.Range("D11") "0.07412"
Set-TextValue $ws.Range("D12") "0.03294"
Set-TextValue $ws.Range("D13") "0.03061"
Set-TextValue $ws.Range("D14") "0.09362"
Set-TextValue $ws.Range("D15") "3.870"
Set-TextValue $ws.Range("D16") "0.001579"
Set-TextValue $ws.Range("D17") "0.04664"

# --- Rows 18-24: coin list shifted by one position with refreshed data ---
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D18") "0.005923"
$ws.Range("E18").Value = "17TigerCashTCH"

$ws.Range("B19").Value = "BitKan"
$ws.Range("C19").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Range("D19") "0.001260"
$ws.Range("E19").Value = "18BitKanKAN"

$ws.Range("B20").Value = "HotbitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Range("D20") "0.004901"
$ws.Range("E20").Value = "19HotbitTokenHTB"

$ws.Range("B21").Value = "NitroEx"
$ws.Range("C21").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Range("D21") "0.00006801"
$ws.Range("E21").Value = "20NitroExNTX"

$ws.Range("B22").Value = "LEO"
$ws.Range("C22").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D22") "3.564"
$ws.Range("E22").Value = "21LEOLEO"

$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D23") "2.144"
$ws.Range("E23").Value = "22BTSETokenBTSE"

$ws.Range("B24").Value = "One"
$ws.Range("C24").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D24") "0.01123"
$ws.Range("E24").Value = "23OneONEBestin24h"

# --- Further column D price updates lower in the sheet ---
Set-TextValue $ws.Range("D40") "0.03971"
Set-TextValue $ws.Range("D41") "0.006183"
Set-TextValue $ws.Range("D42") "0.1071"
Set-TextValue $ws.Range("D43") "0.003001"
Set-TextValue $ws.Range("D44") "0.008727"
Set-TextValue $ws.Range("D45") "0.00005180"
Set-TextValue $ws.Range("D48") "0.002391"
